$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("L3").Value = "stimuli/img_x9w7o.png"
$ws.Range("M3").Value = 92.38888888888889
$ws.Range("N3").Value = 72.94444444444444
$ws.Range("O3").Value = 82.66666666666666
$ws.Range("P3").Value = 36
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = 10
$ws.Range("U3").Value = 10
$ws.Range("V3").Value = 10

# Row 4
$ws.Range("H4").Value = "bedrooms"
$ws.Range("I4").Value = "distractor"
$ws.Range("K4").Value = "f"
$ws.Range("L4").Value = "stimuli/img_jr212.png"
$ws.Range("M4").Value = 59.48571428571429
$ws.Range("N4").Value = 39.08571428571429
$ws.Range("O4").Value = 49.28571428571429
$ws.Range("P4").Value = 35
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 3
$ws.Range("V4").Value = 3

# Row 5
$ws.Range("H5").Value = "bedrooms"
$ws.Range("L5").Value = "stimuli/img_ds88w.png"
$ws.Range("M5").Value = 61.425
$ws.Range("N5").Value = 38.2
$ws.Range("O5").Value = 49.8125
$ws.Range("P5").Value = 40
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 3
$ws.Range("V5").Value = 3

# Row 6
$ws.Range("H6").Value = "living_rooms"
$ws.Range("I6").Value = "target"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_bbs77.png"
$ws.Range("M6").Value = 31.64444444444445
$ws.Range("N6").Value = 21.26666666666667
$ws.Range("O6").Value = 26.45555555555556
$ws.Range("P6").Value = 45
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 2

# Row 7
$ws.Range("H7").Value = "kitchens"
$ws.Range("L7").Value = "stimuli/img_z10c7.png"
$ws.Range("M7").Value = 63.45945945945946
$ws.Range("N7").Value = 37.21621621621622
$ws.Range("O7").Value = 50.33783783783784
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 3
$ws.Range("S7").Value = 3
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 3
$ws.Range("V7").Value = 3

# Row 9
$ws.Range("H9").Value = "kitchens"
$ws.Range("I9").Value = "distractor"
$ws.Range("K9").Value = "f"
$ws.Range("L9").Value = "stimuli/img_1r2ri.png"
$ws.Range("M9").Value = 62.44117647058823
$ws.Range("N9").Value = 40.76470588235294
$ws.Range("O9").Value = 51.60294117647059
$ws.Range("P9").Value = 34
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 3
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 3
$ws.Range("U9").Value = 3
$ws.Range("V9").Value = 3

# Row 10
$ws.Range("L10").Value = "stimuli/img_37hgm.png"
$ws.Range("M10").Value = 70.95454545454545
$ws.Range("N10").Value = 54.77272727272727
$ws.Range("O10").Value = 62.86363636363636
$ws.Range("P10").Value = 44
$ws.Range("Q10").Value = 6
$ws.Range("R10").Value = 6
$ws.Range("S10").Value = 6
$ws.Range("T10").Value = 6
$ws.Range("U10").Value = 6
$ws.Range("V10").Value = 6

# Row 11
$ws.Range("L11").Value = "stimuli/img_rru0v.png"
$ws.Range("M11").Value = 56.45238095238095
$ws.Range("N11").Value = 39.42857142857143
$ws.Range("O11").Value = 47.94047619047619
$ws.Range("P11").Value = 42
$ws.Range("Q11").Value = 4
$ws.Range("R11").Value = 4
$ws.Range("S11").Value = 4
$ws.Range("T11").Value = 4
$ws.Range("U11").Value = 4
$ws.Range("V11").Value = 4

# Row 12
$ws.Range("L12").Value = "stimuli/img_qz292.png"
$ws.Range("M12").Value = 78.26666666666667
$ws.Range("N12").Value = 59.13333333333333
$ws.Range("O12").Value = 68.7
$ws.Range("P12").Value = 45
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 7
$ws.Range("V12").Value = 7

# Row 13
$ws.Range("L13").Value = "stimuli/img_mdpr4.png"
$ws.Range("M13").Value = 74.04255319148936
$ws.Range("N13").Value = 54.70212765957447
$ws.Range("O13").Value = 64.37234042553192
$ws.Range("P13").Value = 47
$ws.Range("Q13").Value = 6
$ws.Range("R13").Value = 6
$ws.Range("S13").Value = 6
$ws.Range("T13").Value = 6
$ws.Range("U13").Value = 6
$ws.Range("V13").Value = 6

# Row 14
$ws.Range("H14").Value = "living_rooms"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_eiu3c.png"
$ws.Range("M14").Value = 65.15909090909091
$ws.Range("N14").Value = 46.22727272727273
$ws.Range("O14").Value = 55.69318181818181
$ws.Range("P14").Value = 44
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 5
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 5
$ws.Range("U14").Value = 5
$ws.Range("V14").Value = 5

# Row 15
$ws.Range("L15").Value = "stimuli/img_bj99b.png"
$ws.Range("M15").Value = 82.79069767441861
$ws.Range("N15").Value = 65.46511627906976
$ws.Range("O15").Value = 74.12790697674419
$ws.Range("Q15").Value = 8
$ws.Range("R15").Value = 8
$ws.Range("S15").Value = 8
$ws.Range("T15").Value = 8
$ws.Range("U15").Value = 8
$ws.Range("V15").Value = 8

# Row 16
$ws.Range("H16").Value = "bedrooms"
$ws.Range("I16").Value = "distractor"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_49h18.png"
$ws.Range("M16").Value = 59.28947368421053
$ws.Range("N16").Value = 32.8421052631579
$ws.Range("O16").Value = 46.06578947368421
$ws.Range("P16").Value = 38
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 3
$ws.Range("T16").Value = 3
$ws.Range("U16").Value = 3
$ws.Range("V16").Value = 3

# Row 17
$ws.Range("H17").Value = "living_rooms"
$ws.Range("I17").Value = "target"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_9oofc.png"
$ws.Range("M17").Value = 82.47619047619048
$ws.Range("N17").Value = 65.5
$ws.Range("O17").Value = 73.98809523809524
$ws.Range("P17").Value = 42
$ws.Range("Q17").Value = 8
$ws.Range("R17").Value = 8
$ws.Range("S17").Value = 8
$ws.Range("T17").Value = 8
$ws.Range("U17").Value = 8
$ws.Range("V17").Value = 8

# Row 19
$ws.Range("H19").Value = "living_rooms"
$ws.Range("I19").Value = "target"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_5nlnv.png"
$ws.Range("M19").Value = 86.1219512195122
$ws.Range("N19").Value = 69.19512195121951
$ws.Range("O19").Value = 77.65853658536585
$ws.Range("P19").Value = 41
$ws.Range("Q19").Value = 9
$ws.Range("R19").Value = 9
$ws.Range("S19").Value = 9
$ws.Range("T19").Value = 9
$ws.Range("U19").Value = 9
$ws.Range("V19").Value = 9

# Row 20
$ws.Range("H20").Value = "kitchens"
$ws.Range("L20").Value = "stimuli/img_2b8fp.png"
$ws.Range("M20").Value = 73.89189189189189
$ws.Range("N20").Value = 51.45945945945946
$ws.Range("O20").Value = 62.67567567567568
$ws.Range("P20").Value = 37
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6
$ws.Range("T20").Value = 6
$ws.Range("U20").Value = 6
$ws.Range("V20").Value = 6

# Row 21
$ws.Range("H21").Value = "bedrooms"
$ws.Range("I21").Value = "distractor"
$ws.Range("K21").Value = "f"
$ws.Range("L21").Value = "stimuli/img_ccn2w.png"
$ws.Range("M21").Value = 65.78723404255319
$ws.Range("N21").Value = 40.31914893617022
$ws.Range("O21").Value = 53.05319148936171
$ws.Range("P21").Value = 47
$ws.Range("Q21").Value = 4
$ws.Range("R21").Value = 4
$ws.Range("S21").Value = 4
$ws.Range("T21").Value = 4
$ws.Range("U21").Value = 4
$ws.Range("V21").Value = 4

# Row 22
$ws.Range("L22").Value = "stimuli/img_2dnfy.png"
$ws.Range("M22").Value = 59.97297297297298
$ws.Range("N22").Value = 37.56756756756756
$ws.Range("O22").Value = 48.77027027027027
$ws.Range("P22").Value = 37
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 3
$ws.Range("S22").Value = 3
$ws.Range("T22").Value = 3
$ws.Range("U22").Value = 3
$ws.Range("V22").Value = 3

# Row 23
$ws.Range("H23").Value = "living_rooms"
$ws.Range("I23").Value = "target"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_jpjeg.png"
$ws.Range("M23").Value = 90.90697674418605
$ws.Range("N23").Value = 74.3953488372093
$ws.Range("O23").Value = 82.65116279069767
$ws.Range("P23").Value = 43
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = 10
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = 10
$ws.Range("U23").Value = 10
$ws.Range("V23").Value = 10

# Row 24
$ws.Range("L24").Value = "stimuli/img_iudc4.png"
$ws.Range("M24").Value = 73.625
$ws.Range("N24").Value = 52.275
$ws.Range("O24").Value = 62.95
$ws.Range("P24").Value = 40
$ws.Range("Q24").Value = 6
$ws.Range("R24").Value = 6
$ws.Range("S24").Value = 6
$ws.Range("T24").Value = 6
$ws.Range("U24").Value = 6
$ws.Range("V24").Value = 6

# Row 25
$ws.Range("H25").Value = "living_rooms"
$ws.Range("I25").Value = "target"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_rg4in.png"
$ws.Range("M25").Value = 49.3695652173913
$ws.Range("N25").Value = 30.21739130434782
$ws.Range("O25").Value = 39.79347826086956
$ws.Range("P25").Value = 46

# Row 26
$ws.Range("H26").Value = "living_rooms"
$ws.Range("I26").Value = "target"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_tbs4n.png"
$ws.Range("M26").Value = 78.95744680851064
$ws.Range("N26").Value = 58.97872340425532
$ws.Range("O26").Value = 68.96808510638297
$ws.Range("P26").Value = 47
$ws.Range("Q26").Value = 7
$ws.Range("R26").Value = 7
$ws.Range("S26").Value = 7
$ws.Range("T26").Value = 7
$ws.Range("U26").Value = 7
$ws.Range("V26").Value = 7

# Row 27
$ws.Range("H27").Value = "kitchens"
$ws.Range("I27").Value = "distractor"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_9mcah.png"
$ws.Range("M27").Value = 59.59375
$ws.Range("N27").Value = 36.75
$ws.Range("O27").Value = 48.171875
$ws.Range("P27").Value = 32
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = 3
$ws.Range("S27").Value = 3
$ws.Range("T27").Value = 3
$ws.Range("U27").Value = 3
$ws.Range("V27").Value = 3

# Row 28
$ws.Range("L28").Value = "stimuli/img_pey7u.png"
$ws.Range("M28").Value = 30.34883720930232
$ws.Range("N28").Value = 20.34883720930232
$ws.Range("O28").Value = 25.34883720930232
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("T28").Value = 2
$ws.Range("U28").Value = 2
$ws.Range("V28").Value = 2

# Row 29
$ws.Range("L29").Value = "stimuli/img_k0ze7.png"
$ws.Range("M29").Value = 58.82142857142857
$ws.Range("N29").Value = 37.46428571428572
$ws.Range("O29").Value = 48.14285714285714
$ws.Range("P29").Value = 28
$ws.Range("Q29").Value = 3
$ws.Range("R29").Value = 3
$ws.Range("S29").Value = 3
$ws.Range("T29").Value = 3
$ws.Range("U29").Value = 3
$ws.Range("V29").Value = 3

# Row 30
$ws.Range("H30").Value = "kitchens"
$ws.Range("I30").Value = "distractor"
$ws.Range("K30").Value = "f"
$ws.Range("L30").Value = "stimuli/img_1r08n.png"
$ws.Range("M30").Value = 72.34285714285714
$ws.Range("N30").Value = 49
$ws.Range("O30").Value = 60.67142857142857
$ws.Range("P30").Value = 35
$ws.Range("Q30").Value = 5
$ws.Range("R30").Value = 5
$ws.Range("S30").Value = 5
$ws.Range("T30").Value = 5
$ws.Range("U30").Value = 5
$ws.Range("V30").Value = 5

# Row 31
$ws.Range("H31").Value = "bedrooms"
$ws.Range("I31").Value = "distractor"
$ws.Range("K31").Value = "f"
$ws.Range("L31").Value = "stimuli/img_gsfx4.png"
$ws.Range("M31").Value = 59.94736842105263
$ws.Range("N31").Value = 33.81578947368421
$ws.Range("O31").Value = 46.88157894736842
$ws.Range("P31").Value = 38
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = 3
$ws.Range("S31").Value = 3
$ws.Range("T31").Value = 3
$ws.Range("U31").Value = 3
$ws.Range("V31").Value = 3

# Row 32
$ws.Range("L32").Value = "stimuli/img_lzz3x.png"
$ws.Range("M32").Value = 18.46341463414634
$ws.Range("N32").Value = 11.92682926829268
$ws.Range("O32").Value = 15.19512195121951
$ws.Range("P32").Value = 41
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = 1
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = 1
$ws.Range("U32").Value = 1
$ws.Range("V32").Value = 1

# Row 33
$ws.Range("H33").Value = "kitchens"
$ws.Range("I33").Value = "distractor"
$ws.Range("K33").Value = "f"
$ws.Range("L33").Value = "stimuli/img_di49w.png"
$ws.Range("M33").Value = 73.03125
$ws.Range("N33").Value = 52.25
$ws.Range("O33").Value = 62.640625
$ws.Range("P33").Value = 32
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = 6
$ws.Range("S33").Value = 6
$ws.Range("T33").Value = 6
$ws.Range("U33").Value = 6
$ws.Range("V33").Value = 6

# Row 34
$ws.Range("H34").Value = "kitchens"
$ws.Range("L34").Value = "stimuli/img_x9ml4.png"
$ws.Range("M34").Value = 58.71428571428572
$ws.Range("N34").Value = 38.46428571428572
$ws.Range("O34").Value = 48.58928571428572
$ws.Range("P34").Value = 28

# Row 35
$ws.Range("H35").Value = "kitchens"
$ws.Range("I35").Value = "distractor"
$ws.Range("K35").Value = "f"
$ws.Range("L35").Value = "stimuli/img_3gm8h.png"
$ws.Range("M35").Value = 65.07894736842105
$ws.Range("N35").Value = 43.92105263157895
$ws.Range("O35").Value = 54.5
$ws.Range("P35").Value = 38
$ws.Range("Q35").Value = 4
$ws.Range("R35").Value = 4
$ws.Range("S35").Value = 4
$ws.Range("T35").Value = 4
$ws.Range("U35").Value = 4
$ws.Range("V35").Value = 4

# Row 36
$ws.Range("H36").Value = "living_rooms"
$ws.Range("I36").Value = "target"
$ws.Range("K36").Value = "j"
$ws.Range("L36").Value = "stimuli/img_5tr4v.png"
$ws.Range("M36").Value = 56.86046511627907
$ws.Range("N36").Value = 39.3953488372093
$ws.Range("O36").Value = 48.12790697674419
$ws.Range("P36").Value = 43
$ws.Range("Q36").Value = 4
$ws.Range("R36").Value = 4
$ws.Range("S36").Value = 4
$ws.Range("T36").Value = 4
$ws.Range("U36").Value = 4
$ws.Range("V36").Value = 4

# Row 37
$ws.Range("H37").Value = "living_rooms"
$ws.Range("I37").Value = "target"
$ws.Range("K37").Value = "j"
$ws.Range("L37").Value = "stimuli/img_il020.png"
$ws.Range("M37").Value = 18.85416666666667
$ws.Range("N37").Value = 16.16666666666667
$ws.Range("O37").Value = 17.51041666666667
$ws.Range("P37").Value = 48
$ws.Range("Q37").Value = 1
$ws.Range("R37").Value = 1
$ws.Range("S37").Value = 1
$ws.Range("T37").Value = 1
$ws.Range("U37").Value = 1
$ws.Range("V37").Value = 1

# Row 38
$ws.Range("H38").Value = "living_rooms"
$ws.Range("I38").Value = "target"
$ws.Range("K38").Value = "j"
$ws.Range("L38").Value = "stimuli/img_4o8l0.png"
$ws.Range("M38").Value = 46.02173913043478
$ws.Range("N38").Value = 31.45652173913043
$ws.Range("O38").Value = 38.73913043478261
$ws.Range("P38").Value = 46

# Row 39
$ws.Range("L39").Value = "stimuli/img_196rk.png"
$ws.Range("M39").Value = 86.53488372093024
$ws.Range("N39").Value = 69.46511627906976
$ws.Range("O39").Value = 78
$ws.Range("P39").Value = 43
$ws.Range("Q39").Value = 9
$ws.Range("R39").Value = 9
$ws.Range("S39").Value = 9
$ws.Range("T39").Value = 9
$ws.Range("U39").Value = 9
$ws.Range("V39").Value = 9

# Row 40
$ws.Range("H40").Value = "bedrooms"
$ws.Range("L40").Value = "stimuli/img_qgbyn.png"
$ws.Range("M40").Value = 65.08108108108108
$ws.Range("N40").Value = 40.10810810810811
$ws.Range("O40").Value = 52.5945945945946
$ws.Range("P40").Value = 37

# Row 41
$ws.Range("H41").Value = "bedrooms"
$ws.Range("I41").Value = "distractor"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_die1d.png"
$ws.Range("M41").Value = 75.42857142857143
$ws.Range("N41").Value = 53.30952380952381
$ws.Range("O41").Value = 64.36904761904762
$ws.Range("P41").Value = 42
$ws.Range("Q41").Value = 6
$ws.Range("R41").Value = 6
$ws.Range("S41").Value = 6
$ws.Range("T41").Value = 6
$ws.Range("U41").Value = 6
$ws.Range("V41").Value = 6

Write-Output "Applied 445 cell updates"
